# Updates the division problems' dividend/divisor text throughout the document.
$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "651÷7="; New = "792÷2=" },
    @{ Old = "775÷9="; New = "921÷5=" },
    @{ Old = "418÷9="; New = "976÷7=" },
    @{ Old = "632÷6="; New = "584÷5=" },
    @{ Old = "185÷7="; New = "198÷5=" },
    @{ Old = "376÷8="; New = "968÷3=" },
    @{ Old = "883÷9="; New = "358÷9=" },
    @{ Old = "863÷3="; New = "500÷8=" },
    @{ Old = "570÷2="; New = "634÷6=" },
    @{ Old = "272÷2="; New = "457÷7=" },
    @{ Old = "263÷8="; New = "410÷8=" },
    @{ Old = "569÷7="; New = "355÷4=" },
    @{ Old = "639÷5="; New = "769÷6=" },
    @{ Old = "280÷2="; New = "422÷9=" },
    @{ Old = "623÷2="; New = "652÷6=" },
    @{ Old = "625÷4="; New = "393÷4=" },
    @{ Old = "552÷3="; New = "478÷3=" },
    @{ Old = "541÷8="; New = "561÷9=" },
    @{ Old = "272÷9="; New = "472÷6=" },
    @{ Old = "885÷9="; New = "663÷9=" },
    @{ Old = "757÷3="; New = "897÷7=" },
    @{ Old = "368÷6="; New = "772÷6=" },
    @{ Old = "742÷4="; New = "874÷6=" },
    @{ Old = "564÷3="; New = "825÷3=" },
    @{ Old = "705÷7="; New = "691÷3=" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.New, 2)
}
